# Generate Report for Archive
#
# The "Status" value flips from "Ready for handoff" to "In Translation"
# everywhere it is reported: the per-language roll-up on the "Overview"
# sheet (columns zh-cn / de-de, row 2) and the "Status" column on each of
# the per-locale detail sheets ("zh-cn", "de-de").
#
# Because the new status text is shorter than the old one, the "Status"
# column(s) that were sized to fit it are narrowed to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# ColumnWidth that lands this engine's pixel-snapped column width as close
# as possible to the narrower "fit the new text" width used by the report
# generator.
$statusColWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $statusColWidth
$overview.Columns.Item(6).ColumnWidth = $statusColWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $statusColWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $statusColWidth
